$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.618.79"
$ws.Range("E2").Value = "  +3.95%  "

$ws.Range("D3").Value = "2.768.87"
$ws.Range("E3").Value = "  +5.02%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'116.32"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.82%  "

$ws.Range("D6").Value = "'333.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.82%  "

$ws.Range("E7").Value = "  +1.96%  "

$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("D9").Value = "'0.575"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.40%  "

$ws.Range("D10").Value = "'41.96"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.17%  "

$ws.Range("D11").Value = "'0.0866"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +6.54%  "

$ws.Range("D12").Value = "'20.29"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.19%  "

$ws.Range("E13").Value = "  +2.30%  "

$ws.Range("D14").Value = "'7.67"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.62%  "

$ws.Range("D15").Value = "3.196.25"
$ws.Range("E15").Value = "  +4.80%  "

$ws.Range("D16").Value = "2.771.08"
$ws.Range("E16").Value = "  +4.97%  "

$ws.Range("D17").Value = "'0.892"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.63%  "

$ws.Range("D18").Value = "51.609.82"
$ws.Range("E18").Value = "  +4.13%  "

$ws.Range("E19").Value = "  +11.18%  "

$ws.Range("D20").Value = "'13.54"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.84%  "

$ws.Range("D21").Value = "'6.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.38%  "

$ws.Range("D22").Value = "0.0₃0979"
$ws.Range("E22").Value = "  +3.32%  "

$ws.Range("D23").Value = "'278.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.07%  "

$ws.Range("D24").Value = "'69.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.13%  "

$ws.Range("D25").Value = "'2.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +6.16%  "

$ws.Range("E26").Value = "  +1.89%  "

$ws.Range("E27").Value = "  +0.08%  "

$ws.Range("D28").Value = "'10.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.71%  "

$ws.Range("D29").Value = "'2.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.09%  "

$ws.Range("E30").Value = "  +1.67%  "

$ws.Range("D31").Value = "'35.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.17%  "

$ws.Range("D32").Value = "'50.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.18%  "

$ws.Range("D33").Value = "'5.58"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.76%  "

$ws.Range("D34").Value = "'0.0822"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.85%  "

$ws.Range("E35").Value = "  -0.10%  "

$ws.Range("D36").Value = "'19.06"
$ws.Range("D36").Style = "Normal"

$ws.Range("D37").Value = "'5.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.49%  "

$ws.Range("D38").Value = "'2.10"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.28%  "

$ws.Range("D39").Value = "'3.25"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.92%  "

$ws.Range("E40").Value = "  +8.87%  "

$ws.Range("E41").Value = "  +0.39%  "

$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'23.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.69%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").Value = "'2.31"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.52%  "

$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").Value = "'0.114"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.88%  "

$ws.Range("D45").Value = "'2.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +15.88%  "

$ws.Range("D46").Value = "2.091.23"
$ws.Range("E46").Value = "  +1.46%  "

$ws.Range("E47").Value = "  +2.79%  "

$ws.Range("E48").Value = "  +3.56%  "

$ws.Range("D49").Value = "'5.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +6.11%  "

$ws.Range("D50").Value = "'8.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.19%  "

$ws.Range("D51").Value = "'60.41"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.10%  "
